$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ROUTE_CONFIG")

# TABLE 1: TRANSPORT MODES - zero out lead time / cost per unit defaults
$ws.Range("B6").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0

# TABLE 2: WAREHOUSE CONFIGURATION - zero out cost/capacity per module defaults
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 0
